$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column headers: "Đơn vị mua" (purchase unit) replaces the old
# "Đơn vị tính" in column F, and "Tỷ lệ quy đổi" (conversion ratio) replaces
# the old "Số lượng quy đổi" in column D. Order of assignment matters so the
# new strings land in the shared-string table in the same order as the
# original authored workbook.
$ws.Range("F1").Value = "Đơn vị mua"
$ws.Range("D1").Value = "Tỷ lệ quy đổi"

# Widen the newly meaningful columns (A, D, E, F) to fit their content
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 9.833333333333332
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666

# Update the active cell selection
$ws.Range("H8").Select()
